# Daily update at 8 AM UTC
# Appends the next day's row of data to the Wins Over Time tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2-98 (header in row 1).
# Find the next empty row right after the last used row in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# New day's serial date value (one day after the previous last row),
# and the day's win counts for each column.
$prevDate = $ws.Cells.Item($lastRow, 1).Value2
$newDate = $prevDate + 1

$ws.Cells.Item($newRow, 1).Value = $newDate
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = 229
$ws.Cells.Item($newRow, 3).Value = 232
$ws.Cells.Item($newRow, 4).Value = 229
